# Fruta / hortaliza, semanal
# Insert a new weekly record as row 73, pushing the existing rows 73-126
# down to rows 74-127 (dimension grows from A1:R126 to A1:R127).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 73 (shifts rows 73..126 down to 74..127).
$ws.Rows("73").Insert()

# Populate the newly inserted row 73 with the new weekly price record.
$ws.Range("A73").Value = 7
$ws.Range("B73").Value = 'Terminal Hortofrutícola Agro Chillán'
$ws.Range("C73").Value = 'Ñuble'
$ws.Range("D73").Value = 44942
$ws.Range("E73").Value = 16
$ws.Range("F73").Value = 100112021
$ws.Range("G73").Value = 'Ají'
$ws.Range("H73").Value = 'Americana (o)'
$ws.Range("I73").Value = 'Primera'
$ws.Range("J73").Value = 60
$ws.Range("K73").Value = 13000
$ws.Range("L73").Value = 13500
$ws.Range("M73").Value = 13250
$ws.Range("N73").Value = '$/caja 15 kilos'
$ws.Range("O73").Value = 'Región del Maule'
$ws.Range("P73").Value = 883
$ws.Range("Q73").Value = 15
$ws.Range("R73").Value = 'Hortaliza'
